$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the existing (previously empty) cells G10 and H10 with their values.
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 5

# I10 and J10 are new cells that need the same "highlighted" style as the
# neighboring cells in column I/J (e.g. I9:J9), so copy that formatting first...
$ws.Range("I9:J9").Copy()
$ws.Range("I10:J10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ...then set their values.
$ws.Range("I10").Value = 5
$ws.Range("J10").Value = 5

# Update the active selection to J10 (was J19).
$ws.Range("J10").Select()
